# Apply "client model, process diff case for fonding name, rating article and etc"
# Updates commodity alternative-name list on Лист1 (sheet 1):
#  - strips trailing non-breaking spaces from several existing names
#  - adds a brand new alternative name "дизтопливо" for дизель (C6)
#  - moves selection to B11
#
# New shared strings must be created in this exact order so they line up
# with the indices used by the target workbook (55 .. 66).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# 55: алюминий
$ws.Range("A12").Value = "алюминий"
# 56: дизтопливо  (new alternative name added next to дизель)
$ws.Range("C6").Value = "дизтопливо"
# 57: дизель
$ws.Range("A6").Value = "дизель"
# 58: нефтяной кокс
$ws.Range("A9").Value = "нефтяной кокс"
# 59: никель
$ws.Range("A13").Value = "никель"
# 60: кобальт
$ws.Range("A19").Value = "кобальт"
# 61: жрс
$ws.Range("A20").Value = "жрс"
# 62: карбамид
$ws.Range("A33").Value = "карбамид"
# reuse existing shared string (43: хлорид калия, no trailing nbsp)
$ws.Range("A35").Value = "хлорид калия"
# 63: аммиачная селитра
$ws.Range("A36").Value = "аммиачная селитра"
# 64: диаммонийфосфат
$ws.Range("A38").Value = "диаммонийфосфат"
# 65: апатитовый концетрат
$ws.Range("A40").Value = "апатитовый концетрат"
# 66: npk-удобрения
$ws.Range("A43").Value = "npk-удобрения"

# Move the active selection to B11 as recorded in the saved workbook view
$ws.Range("B11").Select()
